# Apply the "add tasks in the TODO list" edit to the "short term" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("short term")

# 1) Extend the existing task text in A12 with more detail.
$current = $ws.Range("A12").Value2
$ws.Range("A12").Value2 = $current + " DIFFICILE ! Plutôt mettre un AIC infini si le pic dépasse le min ou le max"

# 2) Add another new task row (row 18) using the sheet's default column style.
#    (Written first so the shared-string table gets this text interned before row 17's.)
$ws.Range("A18").Value2 = "63 regarder le souci avec l'échelle des y qui change d'une page à l'autre dans plotfit2pdf (en résidus si outliers ou non notamment)"

# 3) Add a new task row (row 17), matching the fill/border/wrap style of row 16.
$ws.Range("A17").Value2 = "62. ajouter un etst sur les outliers pour les cas excessifs cf. transcripto rainettes 2018"
$ws.Range("A17").Interior.Color = 49407
$ws.Range("A17").Borders.LineStyle = 1
$ws.Range("A17").WrapText = $true
